$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.598.91"
$ws.Range("E2").Value = "  +3.32%  "

$ws.Range("D3").Value = "1.696.97"
$ws.Range("E3").Value = "  +2.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.89"
$ws.Range("E5").Value = "  +2.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3949"
$ws.Range("E7").Value = "  +1.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4023"
$ws.Range("E8").Value = "  +1.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.528"
$ws.Range("E9").Value = "  +4.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08768"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("E13").Value = "  +7.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.30"
$ws.Range("E14").Value = "  +3.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.095"
$ws.Range("E15").Value = "  +10.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001317"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("D17").Value = "1.692.23"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.90"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07072"
$ws.Range("E19").Value = "  +2.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.71"
$ws.Range("E20").Value = "  +3.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.927"
$ws.Range("E21").Value = "  +4.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.19"
$ws.Range("E23").Value = "  +2.94%  "

$ws.Range("D24").Value = "24.591.16"
$ws.Range("E24").Value = "  +3.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.132"
$ws.Range("E25").Value = "  +10.39%  "

$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.90"
$ws.Range("E27").Value = "  +5.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.79"
$ws.Range("E28").Value = "  +2.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "137.06"
$ws.Range("E29").Value = "  +5.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.196"
$ws.Range("E30").Value = "  +1.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.498"
$ws.Range("E31").Value = "  +10.33%  "

$ws.Range("D32").Value = "1.878.82"
$ws.Range("E32").Value = "  +1.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.090"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08576"
$ws.Range("E34").Value = "  +0.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.185"
$ws.Range("E35").Value = "  +8.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.35"
$ws.Range("E36").Value = "  +7.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2739"
$ws.Range("E37").Value = "  +3.63%  "

$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.50"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09117"
$ws.Range("E40").Value = "  +3.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02743"
$ws.Range("E41").Value = "  +9.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.482"
$ws.Range("E42").Value = "  +2.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7671"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7186"
$ws.Range("E44").Value = "  +2.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.68"
$ws.Range("E45").Value = "  +4.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.547"
$ws.Range("E46").Value = "  +5.53%  "

$ws.Range("E47").Value = "  +2.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9996"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.10"
$ws.Range("E49").Value = "  +1.12%  "

$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.329"
$ws.Range("E50").Value = "  +8.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07989"
$ws.Range("E51").Value = "  +2.47%  "
